$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial numbers (Excel 1900 date system) corresponding to the
# former "YYYYQ4" text labels in column A, rows 2-22 (2004Q4..2024Q4 -> Dec 31 of that year)
$dates = @(38352, 38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657)

$startRow = 2
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
